$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the placeholder "127.0.0.1" IPs in column F (rows 2-6) with the
# real server addresses.
$ws.Range("F2").Value = "192.168.1.113"
$ws.Range("F3").Value = "192.168.1.114"
$ws.Range("F4").Value = "192.168.1.115"
$ws.Range("F5").Value = "192.168.1.116"
$ws.Range("F6").Value = "192.168.1.117"

# Row 7 (the old GameServer_2 row) no longer applies on Linux - clear it out,
# leaving the formatted-but-empty cells behind.
$ws.Range("A7:H7").ClearContents()

# Column F (IP) is now wide enough to show the longer dotted addresses.
$ws.Columns.Item(6).ColumnWidth = 14.318181818181818

# Select the (now blank) row 7, matching the saved selection state.
$ws.Range("A7:XFD7").Select()
